# Added Indian MF 1st Stab
# This script inserts 9 new weekly-tracking columns (B:J) in front of the
# existing data, fills the new header row labels, fills the new body
# cells with the "UN" (unchanged) placeholder used throughout the sheet,
# and records the three new analyst rating actions that were captured
# during those weeks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert 9 new columns before column B, shifting the existing weekly
#    columns (B:V) right to (K:AE).
$ws.Range("B1:J1").EntireColumn.Insert()

# 2. New week headers for the freshly inserted columns.
$ws.Range("B1").Value = "Sep_08"
$ws.Range("C1").Value = "Aug_25"
$ws.Range("D1").Value = "Aug_04"
$ws.Range("E1").Value = "Jul_23"
$ws.Range("F1").Value = "Jul_17"
$ws.Range("G1").Value = "Jul_07"
$ws.Range("H1").Value = "Jun_30"
$ws.Range("I1").Value = "Jun_24"
$ws.Range("J1").Value = "Jun_16"

# 3. Default every new weekly cell (rows 2:33) to "UN" (no rating action
#    that week), matching the rest of the grid.
$ws.Range("B2:J33").Value = "UN"

# 4. Record the three analyst actions that happened during these weeks.
#    Sidoti - Set Price Target, Hold, $83.00 (week of Aug_25)
$ws.Range("C12").Value = "8/9/2019,Set Price Target,Hold,$83.00"
$ws.Range("C12").Interior.Color = 13434828

#    Citigroup - Set Price Target, Hold, $95.00 (week of Jul_07)
$ws.Range("G19").Value = "7/3/2019,Set Price Target,Hold,$95.00"
$ws.Range("G19").Interior.Color = 13434828

#    BidaskClub - Downgrades, Buy -> Hold (week of Aug_25)
$ws.Range("C22").Value = "8/24/2019,Downgrades,Buy -> Hold,"
$ws.Range("C22").Interior.Color = 13408767
